$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# Values are stored as literal text (e.g. '320.21', '3.61%'), so a leading
# apostrophe forces text entry, and resetting Style back to Normal avoids
# leaving a stray quote-prefix style applied to the cell.

$ws.Range("D2").Value = "'320.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'3.61%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.00%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.247"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.48%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07727"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.30%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.709"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.43%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'3.77%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-1.54%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1264"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.64%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1869"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.32%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09232"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.70%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04127"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.70%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.35%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001286"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.06%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005887"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'3.96%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'3.349"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.03%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.345"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.33%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3351"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.430"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'21.36%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1355"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.71%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.20%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04027"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.30%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001271"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.03%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004114"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.46%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001273"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.08%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'4.69%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05334"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.90%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007761"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.02%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1315"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.99%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007038"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.33%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'6.96%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.63%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3470"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.80%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006686"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.69%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'42.72%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004210"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'40.21%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.11%"
$ws.Range("E51").Style = "Normal"

Write-Host "Applied 71 cell updates to Price/Volume columns"
